$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Update values on Tabelle1
$ws1.Range("B5").Value = 62100
$ws1.Range("B6").Value = 69300
$ws1.Range("B7").Value = "01.01.2024"

# Remove data validation on B2 (the dropdown list sourced from Tabelle2)
try {
    $ws1.Range("B2").Validation.Delete()
} catch {
    # ignore if there is nothing to delete / not supported
}

# Move selection to A5
$ws1.Range("A5").Select()
